# Update cryptos list worksheet cell values to match the latest scrape
# (commit: "Updated cryptos list on Fri Apr 21 23:36:20 UTC 2023 with GitHub Actions").
#
# Notes:
#  - Price ("D" column) values such as "0.06500" or "10.20" must stay as TEXT
#    (matching the original inline-string cells) rather than being coerced
#    into numbers, which would silently drop significant trailing zeros or
#    otherwise reformat the value. Assigning through .Formula with a leading
#    single-quote is Excel's standard "force text" prefix and keeps the
#    stored value exactly as typed while remaining a string cell.
#  - Volume/percentage ("E" column) values keep their surrounding spaces and
#    are never number-like, so a plain .Value assignment is safe for them.
#  - Rows 29 and 30 had their contents swapped (LidoDAOToken and
#    InternetComputer(DFINITY) traded places) in addition to getting new
#    numbers, so all four columns (B, C, D, E) are rewritten for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.348.49"
$ws.Range("E2").Value = "  -3.66%  "

$ws.Range("D3").Value = "1.853.58"
$ws.Range("E3").Value = "  -5.08%  "

$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").Formula = "'322.78"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Formula = "'0.4487"
$ws.Range("E7").Value = "  -5.78%  "

$ws.Range("D8").Formula = "'0.3829"
$ws.Range("E8").Value = "  -4.91%  "

$ws.Range("D9").Formula = "'48.35"
$ws.Range("E9").Value = "  -9.73%  "

$ws.Range("D10").Formula = "'0.07850"
$ws.Range("E10").Value = "  -6.61%  "

$ws.Range("D11").Formula = "'1.014"
$ws.Range("E11").Value = "  -3.74%  "

$ws.Range("D12").Formula = "'21.33"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").Value = "1.831.42"
$ws.Range("E13").Value = "  -6.91%  "

$ws.Range("D14").Formula = "'5.856"
$ws.Range("E14").Value = "  -4.54%  "

$ws.Range("D15").Formula = "'7.129"
$ws.Range("E15").Value = "  -5.60%  "

$ws.Range("D16").Formula = "'1.002"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").Formula = "'85.82"
$ws.Range("E17").Value = "  -5.20%  "

$ws.Range("D18").Formula = "'0.00001029"
$ws.Range("E18").Value = "  -3.49%  "

$ws.Range("D19").Formula = "'0.06500"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("E20").Value = "  -7.99%  "

$ws.Range("D21").Formula = "'1.001"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("D22").Formula = "'5.469"
$ws.Range("E22").Value = "  -5.96%  "

$ws.Range("D23").Value = "27.302.38"
$ws.Range("E23").Value = "  -3.86%  "

$ws.Range("D24").Formula = "'10.78"
$ws.Range("E24").Value = "  -5.74%  "

$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "2.033.69"
$ws.Range("E26").Value = "  -7.31%  "

$ws.Range("D27").Formula = "'151.58"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("D28").Formula = "'19.38"
$ws.Range("E28").Value = "  -3.86%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Formula = "'5.548"
$ws.Range("E29").Value = "  -6.00%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Formula = "'2.058"
$ws.Range("E30").Value = "  -4.19%  "

$ws.Range("D31").Formula = "'119.72"
$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").Formula = "'0.09332"
$ws.Range("E32").Value = "  -2.74%  "

$ws.Range("D33").Formula = "'1.475"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").Formula = "'0.9339"
$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("D35").Formula = "'3.600"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("D36").Formula = "'5.255"
$ws.Range("E36").Value = "  -5.94%  "

$ws.Range("E37").Value = "  -4.33%  "

$ws.Range("D38").Formula = "'0.05961"
$ws.Range("E38").Value = "  -3.99%  "

$ws.Range("D39").Formula = "'1.202"
$ws.Range("E39").Value = "  -3.41%  "

$ws.Range("D40").Formula = "'8.273"
$ws.Range("E40").Value = "  -7.60%  "

$ws.Range("D41").Formula = "'1.001"
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("D42").Formula = "'0.5894"
$ws.Range("E42").Value = "  -4.86%  "

$ws.Range("D43").Formula = "'0.1851"
$ws.Range("E43").Value = "  -3.36%  "

$ws.Range("D44").Formula = "'10.20"
$ws.Range("E44").Value = "  -8.08%  "

$ws.Range("D45").Formula = "'1.257"
$ws.Range("E45").Value = "  -6.33%  "

$ws.Range("E46").Value = "  -5.16%  "

$ws.Range("D47").Formula = "'12.15"
$ws.Range("E47").Value = "  -5.97%  "

$ws.Range("D48").Formula = "'1.927"
$ws.Range("E48").Value = "  -6.11%  "

$ws.Range("D49").Formula = "'3.363"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").Formula = "'0.06869"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("D51").Formula = "'108.26"
$ws.Range("E51").Value = "  -1.89%  "
